# Updated cryptos list on Wed Aug 21 11:52:49 UTC 2024 with GitHub Actions
# Refresh price / 1h-volume figures (and the Aave/Bittensor ranking swap at rows 42-43).
# D-column values that look like plain numbers are written with a leading "'"
# (quote-prefix) and the style is immediately reset to "Normal" so the cell keeps
# its original (unstyled) text representation instead of Excel coercing it to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.234.01"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "2.572.63"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'555.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.08%  "
$ws.Range("D6").Value = "'141.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "'0.597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.90%  "
$ws.Range("D9").Value = "2.578.80"
$ws.Range("E9").Value = "  -2.43%  "
$ws.Range("D10").Value = "'6.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("E12").Value = "  +11.16%  "
$ws.Range("D13").Value = "'0.350"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("D14").Value = "3.025.26"
$ws.Range("E14").Value = "  -2.16%  "
$ws.Range("D15").Value = "59.227.96"
$ws.Range("E15").Value = "  -1.94%  "
$ws.Range("D16").Value = "'22.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.97%  "
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "2.580.44"
$ws.Range("E18").Value = "  -2.23%  "
$ws.Range("D19").Value = "'4.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").Value = "'337.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").Value = "'10.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("E22").Value = "  +1.34%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("E24").Value = "  +9.45%  "
$ws.Range("D25").Value = "'62.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.31%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  -2.63%  "
$ws.Range("D28").Value = "'7.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  -3.46%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").Value = "'6.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("D33").Value = "'159.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("D34").Value = "'19.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("E36").Value = "  +1.38%  "
$ws.Range("D37").Value = "'0.893"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("D38").Value = "'37.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("D39").Value = "'0.851"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.82%  "
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("E41").Value = "  +1.10%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "'289.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.27%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'138.30"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.48%  "
$ws.Range("D44").Value = "'0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "'0.0971"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("D46").Value = "'0.592"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.76%  "
$ws.Range("E48").Value = "  -2.59%  "
$ws.Range("D49").Value = "'0.0233"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").Value = "'18.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").Value = "1.937.78"
$ws.Range("E51").Value = "  -0.98%  "
